# Applies the OOXML changes described by the commit:
#   "Chua bai quan ly sinh vien"
#
# Real (non-cosmetic) content changes identified from the diff:
#   1. Slide 10: merge "     AccessModifier " + "DataType " runs into a
#      single run "     AccessModifier DataType " (2nd code line in the
#      "class ClassName{...}" code box).
#   2. Slide 16: table cell text "Constructor được gọi ngầm" ->
#      "Constructor có thể gọi ngầm" (rest of the sentence is kept).
#   3. Slide 18: remove the bullet paragraph "Từ khóa this có thể được
#      truyền như một tham số trong phương Constructor."
#   4. Slide 19: bullet text "Constructors, getters, setters" ->
#      "Constructors"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 10 - merge "AccessModifier " and "DataType " runs together
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(4)
$tr10 = $shp10.TextFrame.TextRange
$para10 = $tr10.Paragraphs(4)
$merged = $para10.Characters(1, 29)
$merged.Text = "     AccessModifier DataType "

# ---------------------------------------------------------------------
# 2) Slide 16 - update the Constructor table cell text
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$shp16 = $s16.Shapes.Item(4)
$tbl16 = $shp16.Table
$cell16 = $tbl16.Cell(4, 1)
$cell16.Shape.TextFrame.TextRange.Text = "Constructor có thể gọi ngầm"

# ---------------------------------------------------------------------
# 3) Slide 18 - remove the redundant "this ... Constructor" bullet
# ---------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$shp18 = $s18.Shapes.Item(4)
$tr18 = $shp18.TextFrame.TextRange
$tr18.Paragraphs(5).Delete()

# ---------------------------------------------------------------------
# 4) Slide 19 - shorten the "Constructors, getters, setters" bullet
# ---------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$shp19 = $s19.Shapes.Item(4)
$tr19 = $shp19.TextFrame.TextRange
$tr19.Paragraphs(3).Text = "Constructors"
